$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5936329588014981
$ws1.Range("C2").Value = 0.5528541226215645
$ws1.Range("D2").Value = 0.9794007490636704
$ws1.Range("E2").Value = 0.7067567567567568
$ws1.Range("F2").Value = 0.8484750162232316
$ws1.Range("G2").Value = 0.9511751538891998
$ws1.Range("H2").Value = 0.6857176422730014
$ws1.Range("I2").Value = 523
$ws1.Range("J2").Value = 423
$ws1.Range("K2").Value = 111
$ws1.Range("L2").Value = 11

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9098360655737705
$ws2.Range("C2").Value = 0.2078651685393259
$ws2.Range("D2").Value = 0.3384146341463415

$ws2.Range("B3").Value = 0.5528541226215645
$ws2.Range("C3").Value = 0.9794007490636704
$ws2.Range("D3").Value = 0.7067567567567568

$ws2.Range("B4").Value = 0.5936329588014981
$ws2.Range("C4").Value = 0.5936329588014981
$ws2.Range("D4").Value = 0.5936329588014981
$ws2.Range("E4").Value = 0.5936329588014981

$ws2.Range("B5").Value = 0.7313450940976676
$ws2.Range("C5").Value = 0.5936329588014981
$ws2.Range("D5").Value = 0.5225856954515491

$ws2.Range("B6").Value = 0.7313450940976676
$ws2.Range("C6").Value = 0.5936329588014981
$ws2.Range("D6").Value = 0.5225856954515491

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 111
$ws3.Range("C2").Value = 423
$ws3.Range("B3").Value = 11
$ws3.Range("C3").Value = 523
